# Over ons.docx - "Daniel Phoeng" bio update
# 1) Finish/clean the "Wie ben ik?" paragraph (drop spell-check run splits,
#    extend the closing sentence).
# 2) Extend the "Waarom heb ik voor deze opleiding gekozen?" answer with a
#    couple more sentences.
# 3) Add a new heading + paragraph "Wat wil ik bereiken met deze opleiding?".

$d = $word.ActiveDocument

# --- 1) "Wie ben ik?" paragraph -------------------------------------------
# Retype the whole paragraph as a single clean run: this keeps the existing
# (unchanged) content intact while naturally dropping the spell-check
# <w:proofErr/> markers that used to split "Daniel Phoeng"/"Codam"/"Coding"/
# "PicNic" into their own runs, and extends the final sentence.
$bioPara = $d.Paragraphs(4)
$bioRange = $bioPara.Range
$bioText = $bioRange.Text
# Strip the trailing paragraph mark before reassigning the text.
$bioTextOnly = $bioText.Substring(0, $bioText.Length - 1)
$newBioTextOnly = $bioTextOnly.Replace("om te leren.", "om hard te leren voor een mooie toekomst!")
$bioTarget = $d.Range($bioRange.Start, $bioRange.Start + $bioTextOnly.Length)
$bioTarget.Text = $newBioTextOnly

# --- 2) "Waarom heb ik voor deze opleiding gekozen?" paragraph ------------
$d.Content.Find.Execute(
    "Zelf heb ik altijd al veel interesse gehad in coderen en de resultaten daarvan. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Zelf heb ik altijd al veel interesse gehad in coderen. Ook al ben ik al gestopt met twee opleidingen die gerelateerd zijn aan coderen, wil ik zeker nog door gaan in deze richting. Ik ben bij de ICT-Academie aangekomen door een recommandatie van een vriend.",
    2) | Out-Null

# --- 3) New heading + paragraph --------------------------------------------
# Locate the (now extended) "Waarom heb ik ..." answer paragraph again and
# create the two new (still empty) paragraphs after it *before* applying any
# bold/underline formatting, so the second (body) paragraph does not inherit
# the heading's character formatting.
$answerPara = $d.Paragraphs(6)
$answerPara.Range.InsertParagraphAfter()
$headingPara = $d.Paragraphs(7)
$headingPara.Range.InsertParagraphAfter()

$headingPara = $d.Paragraphs(7)
$headingRange = $headingPara.Range
$headingRange.Font.Bold = $true
$headingRange.Font.BoldBi = $true
$headingRange.Font.Underline = 1
$headingRange.InsertAfter("Wat wil ik bereiken met deze opleiding?")

$bodyPara = $d.Paragraphs(8)
$bodyPara.Range.InsertAfter("Zelf weet ik nog niet helemaal zeker welke richting ik op wil binnen deze wereld, daarvoor is het fijn dat wij hier nog de verschillende kanten verkennen in de verschillende vakken. Zeker weet ik wel dat ik niet in de richting op wil van hardware en/of security, dat lijkt mij toch net wat minder leuk.")

Write-Output "Done."
